$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the style of H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Populate data rows 2-84 for new columns I and J
# Each entry: row number, I value, J value
$data = @(
    @(2, 9, 9),
    @(3, 5, 7),
    @(4, 7, 7),
    @(5, 8, 8),
    @(6, 8, 8),
    @(7, 9, 9),
    @(8, 7, 7),
    @(9, 10, 10),
    @(10, 7, 7),
    @(11, 7, 7),
    @(12, 10, 10),
    @(13, 8, 8),
    @(14, 8, 8),
    @(15, 8, 8),
    @(16, 8, 9),
    @(17, 7, 7),
    @(18, 8, 8),
    @(19, 8, 8),
    @(20, 10, 10),
    @(21, 7, 8),
    @(22, 7, 7),
    @(23, 7, 8),
    @(24, 7, 8),
    @(25, 8, 8),
    @(26, 6, 7),
    @(27, 6, 6),
    @(28, 8, 8),
    @(29, 7, 7),
    @(30, 8, 8),
    @(31, 7, 7),
    @(32, 7, 7),
    @(33, 8, 8),
    @(34, 8, 8),
    @(35, 9, 9),
    @(36, 6, 6),
    @(37, 7, 7),
    @(38, 7, 7),
    @(39, 6, 6),
    @(40, 7, 7),
    @(41, 9, 9),
    @(42, 7, 7),
    @(43, 8, 8),
    @(44, 10, 10),
    @(45, 10, 10),
    @(46, 8, 8),
    @(47, 8, 8),
    @(48, 7, 7),
    @(49, 8, 8),
    @(50, 7, 7),
    @(51, 6, 7),
    @(52, 8, 8),
    @(53, 10, 10),
    @(54, 8, 8),
    @(55, 9, 9),
    @(56, 8, 8),
    @(57, 8, 8),
    @(58, 7, 8),
    @(59, 8, 8),
    @(60, 8, 8),
    @(61, 8, 8),
    @(62, 9, 9),
    @(63, 8, 8),
    @(64, 8, 8),
    @(65, 8, 8),
    @(66, 8, 8),
    @(67, 8, 8),
    @(68, 8, 8),
    @(69, 8, 8),
    @(70, 8, 8),
    @(71, 11, 11),
    @(72, 8, 8),
    @(73, 7, 8),
    @(74, 8, 8),
    @(75, 7, 7),
    @(76, 10, 10),
    @(77, 8, 8),
    @(78, 8, 8),
    @(79, 5, 5),
    @(80, 7, 7),
    @(81, 4, 4),
    @(82, 4, 4),
    @(83, 7, 7),
    @(84, 4, 4)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}
